$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  2  = -5
  3  = 5
  4  = -3
  6  = -9
  7  = -2
  8  = 1
  9  = 2
  10 = 2
  11 = -6
  12 = 9
  13 = -1
  15 = 5
  16 = -1
  17 = 3
  18 = 2
  20 = -2
  21 = 2
  22 = 1
  23 = 2
  24 = 3
  25 = -3
  26 = -3
  27 = -1
  28 = 2
  29 = -1
  31 = 0
  32 = 8
  33 = 2
  34 = 4
  35 = -2
  36 = 1
  37 = -2
}

foreach ($row in $updates.Keys) {
  $ws.Range("F$row").Value = $updates[$row]
}
